$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.095.63'
$ws.Range("D3").Value = '3.129.41'
$ws.Range("E3").Value = '  +2.13%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = "'578.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.76%  '
$ws.Range("D6").Value = "'179.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.88%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.126.11'
$ws.Range("E8").Value = '  +2.17%  '
$ws.Range("D9").Value = "'0.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.73%  '
$ws.Range("D10").Value = "'6.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.55%  '
$ws.Range("D11").Value = "'0.152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.62%  '
$ws.Range("D12").Value = "'0.472"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.33%  '
$ws.Range("D13").Value = "'0.0000242"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.17%  '
$ws.Range("D14").Value = "'36.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.58%  '
$ws.Range("E15").Value = '  +1.12%  '
$ws.Range("D16").Value = '3.652.88'
$ws.Range("E16").Value = '  +2.18%  '
$ws.Range("D17").Value = '68.064.40'
$ws.Range("E17").Value = '  +2.13%  '
$ws.Range("D18").Value = "'7.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.84%  '
$ws.Range("D19").Value = '3.131.92'
$ws.Range("E19").Value = '  +2.69%  '
$ws.Range("D20").Value = "'16.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.53%  '
$ws.Range("D21").Value = "'488.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.25%  '
$ws.Range("D22").Value = "'0.696"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.46%  '
$ws.Range("D23").Value = "'7.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.20%  '
$ws.Range("D24").Value = "'83.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.55%  '
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = "'12.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.46%  '
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").Value = "'2.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.04%  '
$ws.Range("D27").Value = "'10.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.80%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Value = "'8.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.63%  '
$ws.Range("D30").Value = "'2.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.96%  '
$ws.Range("E31").Value = '  +1.23%  '
$ws.Range("D32").Value = "'28.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.97%  '
$ws.Range("D33").Value = "'0.113"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.90%  '
$ws.Range("D34").Value = '0.0₃0954'
$ws.Range("E34").Value = '  +5.16%  '
$ws.Range("E35").Value = '  +0.20%  '
$ws.Range("D36").Value = "'48.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.13%  '
$ws.Range("B37").Value = 'Mantle'
$ws.Range("C37").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D37").Value = "'0.958"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.05%  '
$ws.Range("B38").Value = 'Filecoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D38").Value = "'5.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.51%  '
$ws.Range("D39").Value = "'0.322"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.45%  '
$ws.Range("D40").Value = "'2.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.51%  '
$ws.Range("D41").Value = "'49.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.26%  '
$ws.Range("E42").Value = '  +2.46%  '
$ws.Range("D43").Value = "'8.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.42%  '
$ws.Range("D44").Value = "'2.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.96%  '
$ws.Range("D45").Value = "'390.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.50%  '
$ws.Range("D46").Value = '2.778.49'
$ws.Range("E46").Value = '  +0.95%  '
$ws.Range("D47").Value = "'27.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +10.93%  '
$ws.Range("D48").Value = "'0.0348"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("D49").Value = "'136.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("D51").Value = "'2.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.18%  '
